$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.220.68"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.548.73"
$ws.Range("E3").Value = "  +4.59%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'569.79"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'148.55"
$ws.Range("E6").Value = "  +4.76%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").Value = "2.545.06"
$ws.Range("E9").Value = "  +4.41%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "'0.356"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "'27.42"
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").Value = "3.004.13"
$ws.Range("E15").Value = "  +4.65%  "
$ws.Range("D16").Value = "63.090.60"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "2.534.96"
$ws.Range("E18").Value = "  +4.38%  "
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").Value = "'336.16"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "'4.31"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'65.38"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("D27").Value = "'1.51"
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'8.41"
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("D30").Value = "'7.23"
$ws.Range("E30").Value = "  +8.27%  "
$ws.Range("D31").Value = "0.0₃0819"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").Value = "'1.86"
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").Value = "'177.52"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E34").Value = "  +6.96%  "
$ws.Range("D35").Value = "'415.08"
$ws.Range("E35").Value = "  +11.04%  "
$ws.Range("D36").Value = "'0.398"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").Value = "'18.87"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D40").Value = "'1.77"
$ws.Range("E40").Value = "  +4.17%  "
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'39.33"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "'152.26"
$ws.Range("E43").Value = "  +2.83%  "
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").Value = "'20.73"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").Value = "'0.609"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'0.0523"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("E49").Value = "  +4.94%  "
$ws.Range("D50").Value = "'18.44"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("D51").Value = "'1.79"
$ws.Range("E51").Value = "  +3.07%  "
